$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before column D. This shifts the existing D:K data
# (and its formatting) one column to the right, into E:L.
$ws.Columns("D").Insert()

# The freshly-inserted column D picked up formatting from column C (its
# left neighbour). Re-copy the correct per-row number formats from column E
# (which now holds what used to be in D) back onto the new column D.
$ws.Range("E7:E35").Copy()
$ws.Range("D7:D35").PasteSpecial(-4122)
$ws.Range("E38:E77").Copy()
$ws.Range("D38:D77").PasteSpecial(-4122)
$ws.Range("E80:E102").Copy()
$ws.Range("D80:D102").PasteSpecial(-4122)

# Populate the new column D with the latest (FY2018) financial figures.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 226600
$ws.Range("D9").Value = 33600
$ws.Range("D10").Value = 193000
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 100
$ws.Range("D15").Value = 77500
$ws.Range("D17").Value = 134500
$ws.Range("D18").Value = 92100
$ws.Range("D20").Value = 92600
$ws.Range("D21").Value = 262300
$ws.Range("D22").Value = 43000
$ws.Range("D23").Value = 141800
$ws.Range("D24").Value = 95800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 46000
$ws.Range("D27").Value = 46000
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -92600
$ws.Range("D33").Value = 46000
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 46000
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 46900
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 37000
$ws.Range("D44").Value = "NA"
$ws.Range("D45").Value = 60800
$ws.Range("D46").Value = 144700
$ws.Range("D47").Value = 12400
$ws.Range("D48").Value = 1924300
$ws.Range("D49").Value = "NA"
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 2200
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 2083600
$ws.Range("D57").Value = 89700
$ws.Range("D58").Value = "NA"
$ws.Range("D59").Value = 72000
$ws.Range("D60").Value = 161700
$ws.Range("D61").Value = 613100
$ws.Range("D62").Value = 111700
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 886600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 101700
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 1197000
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 46000
$ws.Range("D83").Value = 77500
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 67200
$ws.Range("D91").Value = -592700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -706500
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 262100
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -377200
